$d = $word.ActiveDocument
$xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
  <w:pPr><w:rPr><w:b/><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr></w:pPr>
  <w:r><w:rPr><w:b/><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>Quick guide</w:t></w:r>
</w:p>
<w:p>
  <w:pPr><w:rPr><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr></w:pPr>
  <w:r><w:rPr><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">To make your website look similar to the </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:rPr><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>WebNews</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:rPr><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> demo, follow the 3 basic configuration steps.</w:t></w:r>
</w:p>
<w:p>
  <w:pPr><w:rPr><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr></w:pPr>
  <w:r><w:rPr><w:b/><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>1. Install the theme</w:t></w:r>
</w:p>
<w:p>
  <w:pPr><w:rPr><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr></w:pPr>
  <w:r><w:rPr><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>To install it is very simple, go to Appearances&gt; Themes&gt; Add new. Search for the theme or upload the .zip file, if you downloaded the theme manually.</w:t></w:r>
</w:p>
<w:p>
  <w:pPr><w:rPr><w:b/><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr></w:pPr>
  <w:r><w:rPr><w:b/><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">2. Install the recommended </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:rPr><w:b/><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>plugin</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
</w:p>
<w:p>
  <w:pPr><w:rPr><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr></w:pPr>
  <w:r><w:rPr><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Install all recommended </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:rPr><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>plugins</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:rPr><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">, mainly M3 Toolbox. When installing the theme for the first time, a message appears at the top informing you of the necessary </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:rPr><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>plugins</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:rPr><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>.</w:t></w:r>
</w:p>
<w:p>
  <w:pPr><w:rPr><w:b/><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr></w:pPr>
  <w:r><w:rPr><w:b/><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>3. Import the demo content</w:t></w:r>
</w:p>
<w:p>
  <w:pPr><w:rPr><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr></w:pPr>
  <w:r><w:rPr><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">The last step is optional, but this is the easiest way to start styling your site using </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:rPr><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>WebNews</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:rPr><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">, as it will reproduce the demo site on its own site. Demo content is available on the theme author's </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:rPr><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>GitHub</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:rPr><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">, download and import using </w:t></w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r><w:rPr><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>All</w:t></w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r><w:rPr><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:rPr><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>ine</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:rPr><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> one </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:rPr><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>wp</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:rPr><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> migration.</w:t></w:r>
</w:p>
<w:p>
  <w:pPr><w:rPr><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr></w:pPr>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:rPr><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>GitHub</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:rPr><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> link: https://github.com/mariomthree/WebNews-WP-Docs</w:t></w:r>
</w:p>

</w:body>
</w:document>
</pkg:xmlData></pkg:part>
</pkg:package>
'@
$d.Content.InsertXML($xml)
